$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "case with 380 kV done" - update computed line power-flow results (pl_mw)
# for Case_5_47 on Sheet1, columns B,C,D,E,F,J,O across rows 2-25.
$ws.Range("B2").Value = 1.335162643757997
$ws.Range("C2").Value = 0.3007722525749443
$ws.Range("D2").Value = 0.221191879479079
$ws.Range("E2").Value = 0.1743160135657646
$ws.Range("F2").Value = 1.18676455464275
$ws.Range("J2").Value = 0.1819088135304554
$ws.Range("O2").Value = 2.688228375106718
$ws.Range("B3").Value = 1.193101117243145
$ws.Range("C3").Value = 0.2627201377786434
$ws.Range("D3").Value = 0.2147181916054137
$ws.Range("E3").Value = 0.1699470556310629
$ws.Range("F3").Value = 1.191627199631064
$ws.Range("J3").Value = 0.1779213376847437
$ws.Range("O3").Value = 2.715148437092097
$ws.Range("B4").Value = 1.105774502176303
$ws.Range("C4").Value = 0.2392675052066693
$ws.Range("D4").Value = 0.2108020149826046
$ws.Range("E4").Value = 0.1673444571717368
$ws.Range("F4").Value = 1.195574216601997
$ws.Range("J4").Value = 0.1755854790394125
$ws.Range("O4").Value = 2.734165832331286
$ws.Range("B5").Value = 1.070164852877781
$ws.Range("C5").Value = 0.2296886941652474
$ws.Range("D5").Value = 0.209221010957009
$ws.Range("E5").Value = 0.1663040130095723
$ws.Range("F5").Value = 1.197423874964471
$ws.Range("J5").Value = 0.1746618415502255
$ws.Range("O5").Value = 2.742539517129615
$ws.Range("B6").Value = 1.064250543346077
$ws.Range("C6").Value = 0.2280968480349088
$ws.Range("D6").Value = 0.2089593878211531
$ws.Range("E6").Value = 0.166132465087788
$ws.Range("F6").Value = 1.19774556305503
$ws.Range("J6").Value = 0.1745101770417037
$ws.Range("O6").Value = 2.743967593513801
$ws.Range("B7").Value = 1.1052943500178
$ws.Range("C7").Value = 0.239138408993
$ws.Range("D7").Value = 0.210780632653254
$ws.Range("E7").Value = 0.1673303438087004
$ws.Range("F7").Value = 1.195598185755117
$ws.Range("J7").Value = 0.1755729082306487
$ws.Range("O7").Value = 2.734276238828414
$ws.Range("B8").Value = 1.286201707444718
$ws.Range("C8").Value = 0.2876705791923655
$ws.Range("D8").Value = 0.2189476337374145
$ws.Range("E8").Value = 0.1727930170370726
$ws.Range("F8").Value = 1.188241355749263
$ws.Range("J8").Value = 0.1805105650978973
$ws.Range("O8").Value = 2.696992880599879
$ws.Range("B9").Value = 1.640101407823295
$ws.Range("C9").Value = 0.382119492807476
$ws.Range("D9").Value = 0.2354249270776307
$ws.Range("E9").Value = 0.184139216127619
$ws.Range("F9").Value = 1.181467067778982
$ws.Range("J9").Value = 0.1910881103852802
$ws.Range("O9").Value = 2.643705853848758
$ws.Range("B10").Value = 1.899529700330447
$ws.Range("C10").Value = 0.4510502535518981
$ws.Range("D10").Value = 0.2478085709550868
$ws.Range("E10").Value = 0.1928621402399244
$ws.Range("F10").Value = 1.181191195775654
$ws.Range("J10").Value = 0.1994094103582853
$ws.Range("O10").Value = 2.616756535287067
$ws.Range("B11").Value = 2.017413814856525
$ws.Range("C11").Value = 0.4823048200655649
$ws.Range("D11").Value = 0.2535018369541007
$ws.Range("E11").Value = 0.1969146005759512
$ws.Range("F11").Value = 1.182094395936929
$ws.Range("J11").Value = 0.2033154431350539
$ws.Range("O11").Value = 2.607170587895126
$ws.Range("B12").Value = 2.062033205660782
$ws.Range("C12").Value = 0.4941249350550265
$ws.Range("D12").Value = 0.2556662529640619
$ws.Range("E12").Value = 0.1984612832275445
$ws.Range("F12").Value = 1.182584944859897
$ws.Range("J12").Value = 0.2048119639959509
$ws.Range("O12").Value = 2.603927060630809
$ws.Range("B13").Value = 2.05242457513657
$ws.Range("C13").Value = 0.4915799520406949
$ws.Range("D13").Value = 0.2551997309854954
$ws.Range("E13").Value = 0.1981276395424842
$ws.Range("F13").Value = 1.182472681038547
$ws.Range("J13").Value = 0.2044888872727455
$ws.Range("O13").Value = 2.604608391963495
$ws.Range("B14").Value = 2.021085105989869
$ws.Range("C14").Value = 0.483277578234663
$ws.Range("D14").Value = 0.2536797350492037
$ws.Range("E14").Value = 0.1970416046603134
$ws.Range("F14").Value = 1.182131773554389
$ws.Range("J14").Value = 0.2034382140099495
$ws.Range("O14").Value = 2.606895984060515
$ws.Range("B15").Value = 2.001886015826244
$ws.Range("C15").Value = 0.478190123291995
$ws.Range("D15").Value = 0.2527497972755555
$ws.Range("E15").Value = 0.1963779524327336
$ws.Range("F15").Value = 1.181942318267545
$ws.Range("J15").Value = 0.2027969124449527
$ws.Range("O15").Value = 2.608347589271858
$ws.Range("B16").Value = 1.891822893805625
$ws.Range("C16").Value = 0.4490055834085638
$ws.Range("D16").Value = 0.247437698094032
$ws.Range("E16").Value = 0.1925989976332971
$ws.Range("F16").Value = 1.181152920933712
$ws.Range("J16").Value = 0.1991565731298977
$ws.Range("O16").Value = 2.617436969079591
$ws.Range("B17").Value = 1.824267782389541
$ws.Range("C17").Value = 0.4310751692003691
$ws.Range("D17").Value = 0.2441941553559417
$ws.Range("E17").Value = 0.190302318261935
$ws.Range("F17").Value = 1.180932553014244
$ws.Range("J17").Value = 0.196954270974615
$ws.Range("O17").Value = 2.623699182089297
$ws.Range("B18").Value = 1.785399640029539
$ws.Range("C18").Value = 0.4207524670915745
$ws.Range("D18").Value = 0.2423341985498979
$ws.Range("E18").Value = 0.1889892708140408
$ws.Range("F18").Value = 1.180902592998393
$ws.Range("J18").Value = 0.1956989167815664
$ws.Range("O18").Value = 2.627552556475848
$ws.Range("B19").Value = 1.772237513224923
$ws.Range("C19").Value = 0.4172557464386273
$ws.Range("D19").Value = 0.2417054217749381
$ws.Range("E19").Value = 0.1885460603415439
$ws.Range("F19").Value = 1.180909054461367
$ws.Range("J19").Value = 0.195275823750606
$ws.Range("O19").Value = 2.628900379027812
$ws.Range("B20").Value = 1.831460418920585
$ws.Range("C20").Value = 0.4329848904779396
$ws.Range("D20").Value = 0.244538852782
$ws.Range("E20").Value = 0.1905459819535977
$ws.Range("F20").Value = 1.18094598981358
$ws.Range("J20").Value = 0.1971875346156082
$ws.Range("O20").Value = 2.623006513256058
$ws.Range("B21").Value = 2.030290846391722
$ws.Range("C21").Value = 0.4857166068559877
$ws.Range("D21").Value = 0.2541259649036931
$ws.Range("E21").Value = 0.1973602711462235
$ws.Range("F21").Value = 1.182227870585663
$ws.Range("J21").Value = 0.2037463498736827
$ws.Range("O21").Value = 2.606213558173465
$ws.Range("B22").Value = 2.160116125128013
$ws.Range("C22").Value = 0.5200903207177703
$ws.Range("D22").Value = 0.2604411889329725
$ws.Range("E22").Value = 0.201884341647812
$ws.Range("F22").Value = 1.183931659366493
$ws.Range("J22").Value = 0.2081343037887535
$ws.Range("O22").Value = 2.597491932432632
$ws.Range("B23").Value = 2.090837707934952
$ws.Range("C23").Value = 0.5017528050880173
$ws.Range("D23").Value = 0.2570661430368233
$ws.Range("E23").Value = 0.1994633144937055
$ws.Range("F23").Value = 1.182942877077267
$ws.Range("J23").Value = 0.2057830783654708
$ws.Range("O23").Value = 2.60193995995067
$ws.Range("B24").Value = 1.82820872065605
$ws.Range("C24").Value = 0.4321215498936226
$ws.Range("D24").Value = 0.2443830001206777
$ws.Range("E24").Value = 0.1904357987219072
$ws.Range("F24").Value = 1.180939613760813
$ws.Range("J24").Value = 0.1970820425542712
$ws.Range("O24").Value = 2.623318880439456
$ws.Range("B25").Value = 1.544460254482601
$ws.Range("C25").Value = 0.3566480535240544
$ws.Range("D25").Value = 0.2309183347169608
$ws.Range("E25").Value = 0.1810018623999952
$ws.Range("F25").Value = 1.182476715971831
$ws.Range("J25").Value = 0.1881303467705351
$ws.Range("O25").Value = 2.655987105482723
